# "Added new fields in PO screen"
#
# The PO header row (row 1) gains 5 new columns:
#   O1 -> "Deliver To Name"     (new; was "Remarks")
#   P1 -> "Raw Material Code"   (new)
#   Q1 -> "Raw Material Name"   (new)
#   R1 -> "Supplier Load"       (new)
#   S1 -> "Supplier Quantity"   (new)
#   T1 -> "Remarks"             (the former O1 header, pushed right)
#
# giving a used range of A1:T1 instead of A1:O1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header that used to sit in O1 ("Remarks") now lives at the very end,
# in T1. Grab its text (Value2 reads the literal string; Value can come
# back as an opaque COM wrapper in this host) before O1 is overwritten.
$remarksText = $ws.Range("O1").Value2

# Push "Remarks" out to the new last column.
$ws.Range("T1").Value = $remarksText

# Fill in the newly inserted headers.
$ws.Range("O1").Value = "Deliver To Name"
$ws.Range("P1").Value = "Raw Material Code"
$ws.Range("Q1").Value = "Raw Material Name"
$ws.Range("R1").Value = "Supplier Load"
$ws.Range("S1").Value = "Supplier Quantity"

# All header cells share the same bold styling as the rest of row 1.
$ws.Range("O1:T1").Font.Bold = $true

# Size the new columns the way Excel leaves them after typing into
# previously-empty bold header cells (best-fit to content).
$ws.Columns("O").ColumnWidth = 15
$ws.Columns("P").ColumnWidth = 17.1666666666667
$ws.Columns("Q").ColumnWidth = 17.8333333333333
$ws.Columns("R").ColumnWidth = 12.3333333333333
$ws.Columns("S").ColumnWidth = 15.8333333333333
$ws.Columns("T").ColumnWidth = 7.66666666666667

# Leave the selection where it was left after making these edits.
$ws.Range("K6").Select() | Out-Null
